$wb = $excel.ActiveWorkbook

# --- Edit 1: Card11 sheet - disambiguate duplicate "Serviced by" header ---
$wsCard11 = $wb.Worksheets.Item("Card11")
$wsCard11.Range("P1").Value = "Serviced by.1"

# --- Edit 2: Card5 sheet - add new column "Event " ---
$ws = $wb.Worksheets.Item("Card5")

# Header for the new column M (match the style used by the other
# header cells in row 1, e.g. bold font, thin border, centered text)
$ws.Range("M1").Value = "Event "
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 previously had completely empty cells in F:K; the refreshed data
# export fills those with the literal placeholder text "nan" (as used
# throughout the rest of the table for missing values).
$ws.Range("F8").Value = "nan"
$ws.Range("G8").Value = "nan"
$ws.Range("H8").Value = "nan"
$ws.Range("I8").Value = "nan"
$ws.Range("J8").Value = "nan"
$ws.Range("K8").Value = "nan"
